$wb = $excel.ActiveWorkbook

# Worksheets
$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsInclude  = $wb.Worksheets.Item("Include from Health Data Conn")

# Rename the second sheet's tab
$wsInclude.Name = "Include from FFV Schedule Cod"

# Update values on the Metadata sheet
$wsMetadata.Range("B2").Value  = "http://linuxforhealth.org/fhir/cdm/ValueSet/ffv-schedule"
$wsMetadata.Range("B3").Value  = "8.0.0"
$wsMetadata.Range("B5").Value  = "FFV Schedule Value Set"
$wsMetadata.Range("B8").Value  = "2022-11-10T16:00:46+00:00"
$wsMetadata.Range("B9").Value  = "LinuxForHealth Team"
$wsMetadata.Range("B11").Value = "LinuxForHealth valuset defining fee for value initiative schedule codes"

# Update value on the renamed Include sheet (System URI)
$wsInclude.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/ffv-schedule"
